# Refresh crypto Price / Volume(1h) columns with the latest scrape.
# A handful of rows (20/21 and 29/30) also swapped rank order, so
# Coin + Link are rewritten there too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''87.881.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.59%  '

$ws.Range("D3").Value = '''3.321.58'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.88%  '

$ws.Range("D5").Value = '''219.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.60%  '

$ws.Range("D6").Value = '''652.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.41%  '

$ws.Range("E7").Value = '  +22.41%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.605'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("D10").Value = '''3.316.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.91%  '

$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("E12").Value = '  +2.22%  '

$ws.Range("E13").Value = '  +1.60%  '

$ws.Range("D14").Value = '''35.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.39%  '

$ws.Range("D15").Value = '''3.930.79'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.86%  '

$ws.Range("E16").Value = '  +1.73%  '

$ws.Range("D17").Value = '''87.634.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.60%  '

$ws.Range("D18").Value = '''3.320.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.05%  '

$ws.Range("D19").Value = '''14.71'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.77%  '

$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D20").Value = '''3.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.31%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''9.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.85%  '

$ws.Range("D22").Value = '''455.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.74%  '

$ws.Range("D23").Value = '''5.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.41%  '

$ws.Range("D24").Value = '''5.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.07%  '

$ws.Range("E25").Value = '  +12.60%  '

$ws.Range("D26").Value = '''3.492.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.88%  '

$ws.Range("D27").Value = '''78.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("D28").Value = '''0.198'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +39.84%  '

$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '''0.0000127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.22%  '

$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = '''0.998'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.10%  '

$ws.Range("D31").Value = '''614.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.70%  '

$ws.Range("D32").Value = '''9.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.34%  '

$ws.Range("E33").Value = '  +6.42%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("E35").Value = '  +2.53%  '

$ws.Range("D36").Value = '''7.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +21.95%  '

$ws.Range("E37").Value = '  -3.96%  '

$ws.Range("D38").Value = '''23.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("E39").Value = '  +4.79%  '

$ws.Range("D40").Value = '''0.423'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.21%  '

$ws.Range("D41").Value = '''21.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.90%  '

$ws.Range("D42").Value = '''0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").Value = '''3.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = '''159.77'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("D46").Value = '''191.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

$ws.Range("E47").Value = '  +5.39%  '

$ws.Range("D48").Value = '''46.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.81%  '

$ws.Range("D49").Value = '''4.47'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.44%  '

$ws.Range("D50").Value = '''0.782'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.54%  '

$ws.Range("D51").Value = '''0.662'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.89%  '
